$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -19.70265970757035
$ws.Range("C2").Value = 2.043615543081721
$ws.Range("D2").Value = -19.70265970757035
$ws.Range("E2").Value = -19.70265970757035
$ws.Range("F2").Value = -19.70265970757035
$ws.Range("G2").Value = -19.70265970757035
$ws.Range("H2").Value = -19.70265970757035
$ws.Range("I2").Value = -19.70265970757035
$ws.Range("J2").Value = -19.70265970757035
$ws.Range("K2").Value = -19.70265970757035

$ws.Range("B3").Value = -19.70265970757035
$ws.Range("C3").Value = -19.70265970757035
$ws.Range("D3").Value = -19.70265970757035
$ws.Range("E3").Value = -19.70265970757035
$ws.Range("F3").Value = -19.70265970757035
$ws.Range("G3").Value = -19.70265970757035
$ws.Range("H3").Value = -19.70265970757035
$ws.Range("I3").Value = 1.27565204850425
$ws.Range("J3").Value = -19.70265970757035
$ws.Range("K3").Value = -19.70265970757035

$ws.Range("B4").Value = -19.70265970757035
$ws.Range("C4").Value = 1.983104131716404
$ws.Range("D4").Value = 1.840176159036494
$ws.Range("E4").Value = -19.70265970757035
$ws.Range("F4").Value = 3.463419532519206
$ws.Range("G4").Value = -19.70265970757035
$ws.Range("H4").Value = 1.457525022122578
$ws.Range("I4").Value = -19.70265970757035
$ws.Range("J4").Value = 0.9266479921702638
$ws.Range("K4").Value = -19.70265970757035

$ws.Range("B5").Value = -19.70265970757035
$ws.Range("C5").Value = 1.688970407118618
$ws.Range("D5").Value = -19.70265970757035
$ws.Range("E5").Value = -19.70265970757035
$ws.Range("F5").Value = -19.70265970757035
$ws.Range("G5").Value = 2.835915028402222
$ws.Range("H5").Value = -19.70265970757035
$ws.Range("I5").Value = -19.70265970757035
$ws.Range("J5").Value = -19.70265970757035
$ws.Range("K5").Value = -19.70265970757035

$ws.Range("B6").Value = -19.70265970757035
$ws.Range("C6").Value = -19.70265970757035
$ws.Range("D6").Value = -19.70265970757035
$ws.Range("E6").Value = -19.70265970757035
$ws.Range("F6").Value = -19.70265970757035
$ws.Range("G6").Value = -19.70265970757035
$ws.Range("H6").Value = -19.70265970757035
$ws.Range("I6").Value = -19.70265970757035
$ws.Range("J6").Value = -19.70265970757035
$ws.Range("K6").Value = -19.70265970757035

$ws.Range("B7").Value = 2.451002717374493
$ws.Range("C7").Value = -19.70265970757035
$ws.Range("D7").Value = -19.70265970757035
$ws.Range("E7").Value = -19.70265970757035
$ws.Range("F7").Value = -19.70265970757035
$ws.Range("G7").Value = -19.70265970757035
$ws.Range("H7").Value = -19.70265970757035
$ws.Range("I7").Value = -19.70265970757035
$ws.Range("J7").Value = -19.70265970757035
$ws.Range("K7").Value = -19.70265970757035

$ws.Range("B8").Value = -19.70265970757035
$ws.Range("C8").Value = -19.70265970757035
$ws.Range("D8").Value = -19.70265970757035
$ws.Range("E8").Value = -19.70265970757035
$ws.Range("F8").Value = -19.70265970757035
$ws.Range("G8").Value = -19.70265970757035
$ws.Range("H8").Value = -19.70265970757035
$ws.Range("I8").Value = -19.70265970757035
$ws.Range("J8").Value = -19.70265970757035
$ws.Range("K8").Value = -19.70265970757035

$ws.Range("B9").Value = 3.861163231968358
$ws.Range("C9").Value = -19.70265970757035
$ws.Range("D9").Value = -19.70265970757035
$ws.Range("E9").Value = -19.70265970757035
$ws.Range("F9").Value = -19.70265970757035
$ws.Range("G9").Value = -19.70265970757035
$ws.Range("H9").Value = -19.70265970757035
$ws.Range("I9").Value = -19.70265970757035
$ws.Range("J9").Value = -19.70265970757035
$ws.Range("K9").Value = -19.70265970757035

$ws.Range("B10").Value = -19.70265970757035
$ws.Range("C10").Value = -19.70265970757035
$ws.Range("D10").Value = -19.70265970757035
$ws.Range("E10").Value = -19.70265970757035
$ws.Range("F10").Value = -19.70265970757035
$ws.Range("G10").Value = -19.70265970757035
$ws.Range("H10").Value = -19.70265970757035
$ws.Range("I10").Value = 1.752561642013577
$ws.Range("J10").Value = -19.70265970757035
$ws.Range("K10").Value = 2.216841242543552

$ws.Range("B11").Value = -19.70265970757035
$ws.Range("C11").Value = -19.70265970757035
$ws.Range("D11").Value = -19.70265970757035
$ws.Range("E11").Value = -19.70265970757035
$ws.Range("F11").Value = -19.70265970757035
$ws.Range("G11").Value = 2.844014597727286
$ws.Range("H11").Value = -19.70265970757035
$ws.Range("I11").Value = -19.70265970757035
$ws.Range("J11").Value = -19.70265970757035
$ws.Range("K11").Value = 1.961742870544576

$ws.Range("B12").Value = -19.70265970757035
$ws.Range("C12").Value = -19.70265970757035
$ws.Range("D12").Value = -19.70265970757035
$ws.Range("E12").Value = -19.70265970757035
$ws.Range("F12").Value = -19.70265970757035
$ws.Range("G12").Value = -19.70265970757035
$ws.Range("H12").Value = -19.70265970757035
$ws.Range("I12").Value = -19.70265970757035
$ws.Range("J12").Value = -19.70265970757035
$ws.Range("K12").Value = -19.70265970757035

$ws.Range("B13").Value = -19.70265970757035
$ws.Range("C13").Value = -19.70265970757035
$ws.Range("D13").Value = -19.70265970757035
$ws.Range("E13").Value = 4.321926488660741
$ws.Range("F13").Value = -19.70265970757035
$ws.Range("G13").Value = -19.70265970757035
$ws.Range("H13").Value = -19.70265970757035
$ws.Range("I13").Value = -19.70265970757035
$ws.Range("J13").Value = 1.688624269938229
$ws.Range("K13").Value = 1.766265070950215

$ws.Range("B14").Value = -19.70265970757035
$ws.Range("C14").Value = -19.70265970757035
$ws.Range("D14").Value = 1.476071392184416
$ws.Range("E14").Value = -19.70265970757035
$ws.Range("F14").Value = -19.70265970757035
$ws.Range("G14").Value = -19.70265970757035
$ws.Range("H14").Value = -19.70265970757035
$ws.Range("I14").Value = -19.70265970757035
$ws.Range("J14").Value = -19.70265970757035
$ws.Range("K14").Value = 1.968947439457402

$ws.Range("B15").Value = -19.70265970757035
$ws.Range("C15").Value = -19.70265970757035
$ws.Range("D15").Value = 1.758920156524949
$ws.Range("E15").Value = -19.70265970757035
$ws.Range("F15").Value = -19.70265970757035
$ws.Range("G15").Value = -19.70265970757035
$ws.Range("H15").Value = -19.70265970757035
$ws.Range("I15").Value = -19.70265970757035
$ws.Range("J15").Value = -19.70265970757035
$ws.Range("K15").Value = -19.70265970757035

$ws.Range("B16").Value = -19.70265970757035
$ws.Range("C16").Value = -19.70265970757035
$ws.Range("D16").Value = -19.70265970757035
$ws.Range("E16").Value = -19.70265970757035
$ws.Range("F16").Value = -19.70265970757035
$ws.Range("G16").Value = -19.70265970757035
$ws.Range("H16").Value = -19.70265970757035
$ws.Range("I16").Value = -19.70265970757035
$ws.Range("J16").Value = 1.939708023425471
$ws.Range("K16").Value = -19.70265970757035

$ws.Range("B17").Value = -19.70265970757035
$ws.Range("C17").Value = 2.053220926392425
$ws.Range("D17").Value = 1.678095681284394
$ws.Range("E17").Value = -19.70265970757035
$ws.Range("F17").Value = -19.70265970757035
$ws.Range("G17").Value = -19.70265970757035
$ws.Range("H17").Value = 2.079516686095405
$ws.Range("I17").Value = 2.063976359418739
$ws.Range("J17").Value = 2.509176614216512
$ws.Range("K17").Value = -19.70265970757035

$ws.Range("B18").Value = -19.70265970757035
$ws.Range("C18").Value = -19.70265970757035
$ws.Range("D18").Value = -19.70265970757035
$ws.Range("E18").Value = -19.70265970757035
$ws.Range("F18").Value = -19.70265970757035
$ws.Range("G18").Value = -19.70265970757035
$ws.Range("H18").Value = 1.99802933741002
$ws.Range("I18").Value = 2.050627647228624
$ws.Range("J18").Value = 2.418569299797367
$ws.Range("K18").Value = -19.70265970757035

$ws.Range("B19").Value = -19.70265970757035
$ws.Range("C19").Value = -19.70265970757035
$ws.Range("D19").Value = 1.993484745983819
$ws.Range("E19").Value = -19.70265970757035
$ws.Range("F19").Value = -19.70265970757035
$ws.Range("G19").Value = -19.70265970757035
$ws.Range("H19").Value = 1.607676263260895
$ws.Range("I19").Value = 1.823113894339387
$ws.Range("J19").Value = -19.70265970757035
$ws.Range("K19").Value = -19.70265970757035

$ws.Range("B20").Value = -19.70265970757035
$ws.Range("C20").Value = 1.025640330907295
$ws.Range("D20").Value = 1.618740510395683
$ws.Range("E20").Value = -19.70265970757035
$ws.Range("F20").Value = 3.165032473689277
$ws.Range("G20").Value = -19.70265970757035
$ws.Range("H20").Value = 1.680070377000127
$ws.Range("I20").Value = 1.230610958784366
$ws.Range("J20").Value = -19.70265970757035
$ws.Range("K20").Value = 2.049410951669235

$ws.Range("B21").Value = -19.70265970757035
$ws.Range("C21").Value = 1.329604668793474
$ws.Range("D21").Value = -19.70265970757035
$ws.Range("E21").Value = -19.70265970757035
$ws.Range("F21").Value = -19.70265970757035
$ws.Range("G21").Value = 2.505858364296701
$ws.Range("H21").Value = 1.475281071206393
$ws.Range("I21").Value = -19.70265970757035
$ws.Range("J21").Value = -19.70265970757035
$ws.Range("K21").Value = -19.70265970757035

